$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data block (rows 2-9) is duplicated twice more, appended
# directly below itself, growing the table from 9 rows to 25 rows.
$sourceRange = $ws.Range("A2:E9")
$sourceRange.Copy()

$ws.Range("A10").PasteSpecial()
$ws.Range("A18").PasteSpecial()

$excel.CutCopyMode = $false

$ws.Range("I23").Select()
